$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 8-10 (no longer present in the updated data)
$ws.Range("A8:T10").EntireRow.Delete() | Out-Null

# Updated data for rows 2-7 (D column cluster label + numeric columns E..T)
# Cluster string ids used via shared strings: ECs, FAPs, MuSCs already present in B/C/A columns
$data = @{
    2 = @{ D = "FAPs";  E = 3; F = 1; G = 26.127733; H = 78.383199; I = 0.2666992864894373; J = 0.2666992864894374; K = 3; L = 1; M = 0.1124773333333333; N = 0.337432; O = 0.7871437602495106; P = 0.7871437602495107; Q = 2.938777733885333; R = 26.448999604968; S = 0.2099306792231572; T = 0.2099306792231573 }
    3 = @{ D = "MuSCs"; E = 3; F = 1; G = 26.127733; H = 78.383199; I = 0.2666992864894373; J = 0.2666992864894374; K = 2; L = 0.6666666666666666; M = 0.03041566666666666; N = 0.09124699999999999; O = 0.2128562397504893; P = 0.2128562397504893; Q = 0.7946924176836667; R = 7.152231759153; S = 0.05676860726628009; T = 0.05676860726628011 }
    4 = @{ A = "FAPs";  D = "FAPs";  E = 3; F = 1; G = 23.90796933333333; H = 71.72390799999999; I = 0.2440410104700376; J = 0.2440410104700377; K = 3; L = 1; M = 0.1124773333333333; N = 0.337432; O = 0.7871437602495106; P = 0.7871437602495107; Q = 2.689104636028444; R = 24.201941724256; S = 0.1920953586364756; T = 0.1920953586364756 }
    5 = @{ D = "MuSCs"; E = 3; F = 1; G = 23.90796933333333; H = 71.72390799999999; I = 0.2440410104700376; J = 0.2440410104700377; K = 2; L = 0.6666666666666666; M = 0.03041566666666666; N = 0.09124699999999999; O = 0.2128562397504893; P = 0.2128562397504893; Q = 0.7271768259195553; R = 6.544591433275999; S = 0.05194565183356199; T = 0.05194565183356201 }
    6 = @{ A = "MuSCs"; D = "FAPs";  E = 3; F = 1; G = 47.93131266666666; H = 143.793938; I = 0.489259703040525; J = 0.4892597030405251; K = 3; L = 1; M = 0.1124773333333333; N = 0.337432; O = 0.7871437602495106; P = 0.7871437602495107; Q = 5.391186231912887; R = 48.52067608721599; S = 0.3851177223898778; T = 0.3851177223898779 }
    7 = @{ A = "MuSCs"; D = "MuSCs"; E = 3; F = 1; G = 47.93131266666666; H = 143.793938; I = 0.489259703040525; J = 0.4892597030405251; K = 2; L = 0.6666666666666666; M = 0.03041566666666666; N = 0.09124699999999999; O = 0.2128562397504893; P = 0.2128562397504893; Q = 1.457862828965111; R = 13.120765460686; S = 0.1041419806506472; T = 0.1041419806506472 }
}

foreach ($r in $data.Keys) {
    $row = $data[$r]
    foreach ($col in $row.Keys) {
        $ws.Range("$col$r").Value = $row[$col]
    }
}
